$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Price column to text so values like "226.89" are not
# auto-converted to numbers by Excel's type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '34.025.94'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '1.783.81'
$ws.Range("E3").Value = '  -2.00%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '226.89'
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").Value = '0.551'
$ws.Range("E6").Value = '  +1.24%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '31.22'
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("D9").Value = '46.20'
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("D10").Value = '0.280'
$ws.Range("E10").Value = '  -0.37%  '
$ws.Range("E11").Value = '  -2.36%  '
$ws.Range("E12").Value = '  -0.23%  '
$ws.Range("D13").Value = '2.038.71'
$ws.Range("E13").Value = '  -2.22%  '
$ws.Range("D14").Value = '11.36'
$ws.Range("E14").Value = '  +10.93%  '
$ws.Range("D15").Value = '1.786.95'
$ws.Range("E15").Value = '  -1.70%  '
$ws.Range("E16").Value = '  -1.67%  '
$ws.Range("D17").Value = '34.018.90'
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("E18").Value = '  -2.76%  '
$ws.Range("D19").Value = '69.34'
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").Value = '252.48'
$ws.Range("E20").Value = '  -2.21%  '
$ws.Range("D21").Value = '0.0₃0740'
$ws.Range("E21").Value = '  -1.12%  '
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("D23").Value = '10.42'
$ws.Range("E23").Value = '  -0.61%  '
$ws.Range("D24").Value = '4.23'
$ws.Range("E24").Value = '  -3.21%  '
$ws.Range("E25").Value = '  -2.49%  '
$ws.Range("D26").Value = '156.98'
$ws.Range("E26").Value = '  -2.30%  '
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D28").Value = '7.01'
$ws.Range("E28").Value = '  -1.97%  '
$ws.Range("E29").Value = '  -1.91%  '
$ws.Range("E30").Value = '  -0.19%  '
$ws.Range("E31").Value = '  -1.99%  '
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").Value = '1.20'
$ws.Range("E33").Value = '  -1.08%  '
$ws.Range("E34").Value = '  +0.90%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = '1.449.97'
$ws.Range("E36").Value = '  -7.82%  '
$ws.Range("E37").Value = '  -0.39%  '
$ws.Range("D38").Value = '0.627'
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("D39").Value = '0.0187'
$ws.Range("E39").Value = '  -1.08%  '
$ws.Range("D40").Value = '83.27'
$ws.Range("E40").Value = '  -1.80%  '
$ws.Range("E41").Value = '  -0.76%  '
$ws.Range("E42").Value = '  +0.57%  '
$ws.Range("D43").Value = '0.897'
$ws.Range("E43").Value = '  -1.87%  '
$ws.Range("E44").Value = '  -2.74%  '
$ws.Range("E45").Value = '  -2.32%  '
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("D47").Value = '1.939.68'
$ws.Range("E47").Value = '  -2.11%  '
$ws.Range("D48").Value = '5.74'
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").Value = '11.86'
$ws.Range("E50").Value = '  +5.80%  '
$ws.Range("D51").Value = '50.97'
$ws.Range("E51").Value = '  -3.44%  '

# Restore original (default) cell style now that text values are set.
$ws.Range("D2:D51").Style = "Normal"
